# "Timing Results for BCGS" - fill in newly-available timing numbers for the
# full_hypre (3) and full_hypre (4) sheets, and move the active selection from
# full_hypre (2) (D11) to full_hypre (4).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# full_hypre (3): rows 4-8 get complete timing data (C..K); rows 2-3 only get
# empty/blank touched cells in D..K (their source rows had no values yet).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("full_hypre (3)")

$ws3.Range("D2:K3").NumberFormat = "General"

$ws3.Range("C4").Value = 26
$ws3.Range("D4").Value = [double]"1.7253507721289E-10"
$ws3.Range("E4").Value = [double]"0.00023017137239509"
$ws3.Range("F4").Value = [double]"0.2214224821"
$ws3.Range("G4").Value = [double]"2.696783343"
$ws3.Range("H4").Value = [double]"446.9183317123"
$ws3.Range("I4").Value = [double]"449.6151316394"
$ws3.Range("J4").Value = [double]"47.3735166304"
$ws3.Range("K4").Value = [double]"47.3735270082"

$ws3.Range("C5").Value = 26
$ws3.Range("D5").Value = [double]"9.506452256118E-11"
$ws3.Range("E5").Value = [double]"0.00023017137260259"
$ws3.Range("F5").Value = [double]"0.1196586295"
$ws3.Range("G5").Value = [double]"1.5033511788"
$ws3.Range("H5").Value = [double]"181.5855268028"
$ws3.Range("I5").Value = [double]"183.0888965869"
$ws3.Range("J5").Value = [double]"26.7883942205"
$ws3.Range("K5").Value = [double]"26.7884045252"

$ws3.Range("C6").Value = 25
$ws3.Range("D6").Value = [double]"4.1202575554487E-10"
$ws3.Range("E6").Value = [double]"0.00023017137259984"
$ws3.Range("F6").Value = [double]"0.0616104063"
$ws3.Range("G6").Value = [double]"0.7892406357"
$ws3.Range("H6").Value = [double]"54.2893862175"
$ws3.Range("I6").Value = [double]"55.0786460205"
$ws3.Range("J6").Value = [double]"14.1369411747"
$ws3.Range("K6").Value = [double]"14.1369528537"

$ws3.Range("C7").Value = 25
$ws3.Range("D7").Value = [double]"2.1046551056235E-10"
$ws3.Range("E7").Value = [double]"0.0002301713723063"
$ws3.Range("F7").Value = [double]"0.0325132696"
$ws3.Range("G7").Value = [double]"0.4025836837"
$ws3.Range("H7").Value = [double]"16.1595630834"
$ws3.Range("I7").Value = [double]"16.5621819183"
$ws3.Range("J7").Value = [double]"7.8431197993"
$ws3.Range("K7").Value = [double]"7.843144716"

$ws3.Range("C8").Value = 25
$ws3.Range("D8").Value = [double]"4.8162963120024E-10"
$ws3.Range("E8").Value = [double]"0.00023017137249261"
$ws3.Range("F8").Value = [double]"0.0190281091"
$ws3.Range("G8").Value = [double]"0.2109954786"
$ws3.Range("H8").Value = [double]"5.1019151299"
$ws3.Range("I8").Value = [double]"5.3129964661"
$ws3.Range("J8").Value = [double]"5.6299329695"
$ws3.Range("K8").Value = [double]"5.6299969828"

# ---------------------------------------------------------------------------
# full_hypre (4): rows 6-8 get complete timing data (C..K).
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("full_hypre (4)")

$ws5.Range("C6").Value = 27
$ws5.Range("D6").Value = [double]"4.0835231674808E-10"
$ws5.Range("E6").Value = [double]"5.7477243126367E-05"
$ws5.Range("F6").Value = [double]"0.2474672632"
$ws5.Range("G6").Value = [double]"3.1432230861"
$ws5.Range("H6").Value = [double]"758.8059399162"
$ws5.Range("I6").Value = [double]"761.9491857614"
$ws5.Range("J6").Value = [double]"59.6765665637"
$ws5.Range("K6").Value = [double]"59.6765805297"

$ws5.Range("C7").Value = 28
$ws5.Range("D7").Value = [double]"4.4031378066206E-10"
$ws5.Range("E7").Value = [double]"5.7477243317198E-05"
$ws5.Range("F7").Value = [double]"0.1324929758"
$ws5.Range("G7").Value = [double]"1.6267713022"
$ws5.Range("H7").Value = [double]"223.1720150415"
$ws5.Range("I7").Value = [double]"224.7988406424"
$ws5.Range("J7").Value = [double]"42.2347568727"
$ws5.Range("K7").Value = [double]"42.2347812589"

$ws5.Range("C8").Value = 26
$ws5.Range("D8").Value = [double]"1.4853340375077E-09"
$ws5.Range("E8").Value = [double]"5.7477243111602E-05"
$ws5.Range("F8").Value = [double]"0.0689320307"
$ws5.Range("G8").Value = [double]"0.9686685738"
$ws5.Range("H8").Value = [double]"63.1584382067"
$ws5.Range("I8").Value = [double]"64.1272050631"
$ws5.Range("J8").Value = [double]"21.7812048834"
$ws5.Range("K8").Value = [double]"21.7812723746"

# ---------------------------------------------------------------------------
# Switch the active sheet/selection: the user was on full_hypre (2) with D11
# selected, and moves over to full_hypre (4) (which keeps its default A1
# selection) to review the newly added BCGS timings.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("full_hypre (2)")
$ws1.Activate()
$ws1.Range("D11").Select() | Out-Null

$ws5.Activate()
